$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B22").Value = 7202436
$ws.Range("E22").Value = "Lucko"
$ws.Range("F22").Value = "NK Tondach"
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 1
$ws.Range("J22").Value = 1.615
$ws.Range("K22").Value = 3.75
$ws.Range("L22").Value = 4.333
$ws.Range("M22").Value = 1.533
$ws.Range("N22").Value = 4
$ws.Range("O22").Value = 5
$ws.Range("P22").Value = -0.75
$ws.Range("Q22").Value = 1.7
$ws.Range("R22").Value = 2.1
$ws.Range("V22").Value = 0.5329999999999999
$ws.Range("Y22").Value = 0.35
$ws.Range("Z22").Value = -0.5
$ws.Range("AA22").Value = 0.4125
$ws.Range("AB22").Value = -0.5

$ws.Range("B23").Value = 7202435
$ws.Range("E23").Value = "NK Udarnik Kurilovec"
$ws.Range("F23").Value = "NK Mladost Petrinja"
$ws.Range("J23").Value = 2
$ws.Range("K23").Value = 3.4
$ws.Range("L23").Value = 3.1
$ws.Range("M23").Value = 2
$ws.Range("N23").Value = 3.4
$ws.Range("O23").Value = 3.1
$ws.Range("P23").Value = -0.25
$ws.Range("Q23").Value = 1.8
$ws.Range("R23").Value = 2
$ws.Range("S23").Value = 2.75
$ws.Range("T23").Value = 1.825
$ws.Range("U23").Value = 1.975
$ws.Range("V23").Value = 1
$ws.Range("Y23").Value = 0.8
$ws.Range("AA23").Value = 0.825

$ws.Range("B24").Value = 7202437
$ws.Range("E24").Value = "NK Maksimir"
$ws.Range("F24").Value = "Sava Strmec"
$ws.Range("G24").Value = 6
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 1.4
$ws.Range("K24").Value = 4.333
$ws.Range("L24").Value = 6
$ws.Range("M24").Value = 1.4
$ws.Range("N24").Value = 4.333
$ws.Range("O24").Value = 6
$ws.Range("P24").Value = -1.25
$ws.Range("Q24").Value = 1.85
$ws.Range("R24").Value = 1.95
$ws.Range("S24").Value = 3
$ws.Range("T24").Value = 1.8
$ws.Range("U24").Value = 2
$ws.Range("V24").Value = 0.3999999999999999
$ws.Range("Y24").Value = 0.8500000000000001
$ws.Range("Z24").Value = -1
$ws.Range("AA24").Value = 0.8
$ws.Range("AB24").Value = -1

$ws.Range("B33").Value = 7291473
$ws.Range("E33").Value = "NK Maksimir"
$ws.Range("F33").Value = "NK Mladost Petrinja"
$ws.Range("G33").Value = 5
$ws.Range("J33").Value = 1.25
$ws.Range("K33").Value = 6
$ws.Range("L33").Value = 7
$ws.Range("M33").Value = 1.25
$ws.Range("N33").Value = 6
$ws.Range("O33").Value = 7.5
$ws.Range("P33").Value = -1.75
$ws.Range("Q33").Value = 1.9
$ws.Range("R33").Value = 1.9
$ws.Range("S33").Value = 3
$ws.Range("T33").Value = 1.975
$ws.Range("U33").Value = 1.825
$ws.Range("V33").Value = 0.25
$ws.Range("Y33").Value = 0.8999999999999999
$ws.Range("AA33").Value = 0.9750000000000001

$ws.Range("B34").Value = 7291472
$ws.Range("E34").Value = "NK Lukavec"
$ws.Range("F34").Value = "Sava Strmec"
$ws.Range("G34").Value = 3
$ws.Range("J34").Value = 2.2
$ws.Range("K34").Value = 3.6
$ws.Range("L34").Value = 2.6
$ws.Range("M34").Value = 2.2
$ws.Range("N34").Value = 3.6
$ws.Range("O34").Value = 2.625
$ws.Range("P34").Value = -0.25
$ws.Range("Q34").Value = 2
$ws.Range("R34").Value = 1.8
$ws.Range("S34").Value = 2.5
$ws.Range("T34").Value = 1.8
$ws.Range("U34").Value = 2
$ws.Range("V34").Value = 1.2
$ws.Range("Y34").Value = 1
$ws.Range("AA34").Value = 0.8

$ws.Range("B124").Value = 8152457
$ws.Range("E124").Value = "NK Tomislav Cerna"
$ws.Range("F124").Value = "NK Tomislav Donji Andrijevci"
$ws.Range("H124").Value = 0
$ws.Range("I124").Value = "H"
$ws.Range("K124").Value = 4
$ws.Range("L124").Value = 2.25
$ws.Range("M124").Value = 1.7
$ws.Range("N124").Value = 3.75
$ws.Range("O124").Value = 3.8
$ws.Range("P124").Value = -0.75
$ws.Range("T124").Value = 1.85
$ws.Range("U124").Value = 1.95
$ws.Range("V124").Value = 0.7
$ws.Range("X124").Value = -1
$ws.Range("Y124").Value = 0.4625
$ws.Range("Z124").Value = -0.5
$ws.Range("AA124").Value = -1
$ws.Range("AB124").Value = 0.95

$ws.Range("B125").Value = 8152456
$ws.Range("E125").Value = "Sava Strmec"
$ws.Range("F125").Value = "NK Samobor"
$ws.Range("H125").Value = 3
$ws.Range("I125").Value = "A"
$ws.Range("K125").Value = 3.6
$ws.Range("L125").Value = 2.4
$ws.Range("M125").Value = 2.1
$ws.Range("N125").Value = 3.6
$ws.Range("O125").Value = 2.75
$ws.Range("P125").Value = -0.25
$ws.Range("T125").Value = 1.775
$ws.Range("U125").Value = 2.025
$ws.Range("V125").Value = -1
$ws.Range("X125").Value = 1.75
$ws.Range("Y125").Value = -1
$ws.Range("Z125").Value = 0.875
$ws.Range("AA125").Value = 0.7749999999999999
$ws.Range("AB125").Value = -1
